$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-30 Saturday" "2024-03-31 Sunday"

Replace-Text "13×26=" "65×23="
Replace-Text "97×14=" "88×26="
Replace-Text "66×15=" "69×29="
Replace-Text "74×17=" "63×95="
Replace-Text "92×50=" "81×53="

Replace-Text "30×12=" "19×53="
Replace-Text "43×46=" "80×75="
Replace-Text "44×77=" "36×61="
Replace-Text "77×30=" "63×67="
Replace-Text "42×20=" "35×95="

Replace-Text "55×33=" "51×19="
Replace-Text "90×45=" "58×16="
Replace-Text "41×41=" "55×53="
Replace-Text "89×69=" "11×89="
Replace-Text "95×13=" "39×77="

Replace-Text "56×56=" "73×37="
Replace-Text "31×41=" "16×56="
Replace-Text "69×90=" "44×97="
Replace-Text "15×36=" "66×72="
Replace-Text "54×59=" "91×13="

Replace-Text "63×84=" "13×39="
Replace-Text "40×57=" "84×75="
Replace-Text "61×40=" "40×18="
Replace-Text "96×44=" "65×97="
Replace-Text "71×78=" "26×23="
